# Add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row 2 holding the 2022-Q3 totals,
#    pushing the existing quarters down by one row.
# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e. before
#    the existing "2022-Q2" tab) and fill it with the per-fund breakdown.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" sheet.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push rows 2..N down by inserting a fresh row at position 2.
$summary.Rows.Item(2).Insert()

# Row-insert in Excel copies the format of the row above into the new row;
# reset the data cells (B:D) back to the unstyled look the other data rows
# use, then restyle the index cell (A) like the rest of column A.
$summary.Range("B2:D2").ClearFormats()

$summary.Range("A2").Value = 0
$summary.Range("A2").Font.Bold = $true
$summary.Range("A2").HorizontalAlignment = -4108
$summary.Range("A2").VerticalAlignment = -4160
$summary.Range("A2").Borders.LineStyle = 1

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 16
$summary.Range("D2").Value = 2.98

# The row-insert left the old index values (0,1,2,3) in place on rows 3..6;
# bump them by one (1,2,3,4) to keep the running index consistent.
for ($r = 3; $r -le 6; $r++) {
    $summary.Cells.Item($r, 1).Value = $r - 2
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet before the existing "2022-Q2" tab.
# ---------------------------------------------------------------------------
$before = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($before)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $cell = $q3.Cells.Item(1, $c + 2)
    $cell.Value = $headers[$c]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

$rows = @(
  @("0","206009","鹏华新兴产业混合","41.17","89.21","4.41","1.8156","8"),
  @("1","003713","英大睿盛灵活配置混合A","2.83","93.65","9.49","0.2686","1"),
  @("2","003714","英大睿盛灵活配置混合C","2.19","93.65","9.49","0.2078","1"),
  @("3","008132","鹏华价值驱动混合","4.47","92.01","4.16","0.1860","3"),
  @("4","519655","银河现代服务主题灵活配置混合","3.37","85.06","4.50","0.1516","7"),
  @("5","001678","英大国企改革主题股票","1.55","93.30","5.50","0.0852","5"),
  @("6","007832","博道伍佰智航股票C","6.02","88.19","1.01","0.0608","6"),
  @("7","160323","华夏磐泰混合（LOF）A","6.20","28.65","0.73","0.0453","1"),
  @("8","001607","英大策略优选混合A","0.57","91.98","5.78","0.0329","7"),
  @("9","007831","博道伍佰智航股票A","3.13","88.19","1.01","0.0316","6"),
  @("10","012522","英大稳固增强核心一年持有混合C","1.24","27.71","2.33","0.0289","1"),
  @("11","013360","华夏磐泰混合（LOF）C","3.70","28.65","0.73","0.0270","1"),
  @("12","012521","英大稳固增强核心一年持有混合A","0.75","27.71","2.33","0.0175","1"),
  @("13","003447","英大睿鑫灵活配置混合C","0.21","92.71","7.76","0.0163","6"),
  @("14","003446","英大睿鑫灵活配置混合A","0.07","92.71","7.76","0.0054","6"),
  @("15","001608","英大策略优选混合C","0.02","91.98","5.78","0.0012","7")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]

    $idxCell = $q3.Cells.Item($r, 1)
    $idxCell.Value = [int]$row[0]
    $idxCell.Font.Bold = $true
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160
    $idxCell.Borders.LineStyle = 1

    $q3.Cells.Item($r, 2).Value = "'" + $row[1]
    $q3.Cells.Item($r, 3).Value = "'" + $row[2]
    $q3.Cells.Item($r, 4).Value = "'" + $row[3]
    $q3.Cells.Item($r, 5).Value = "'" + $row[4]
    $q3.Cells.Item($r, 6).Value = "'" + $row[5]
    $q3.Cells.Item($r, 7).Value = "'" + $row[6]
    $q3.Cells.Item($r, 8).Value = [int]$row[7]
}

# Keep the originally-active tab (the last sheet) selected, same as before
# the edit - inserting/populating sheets shouldn't change the user's view.
$wb.Worksheets.Item($wb.Worksheets.Count).Activate()
